$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.191581249237061
$ws.Range("B1").Value = 2.490107297897339
$ws.Range("C1").Value = 2.140068292617798
$ws.Range("D1").Value = 2.27028751373291
$ws.Range("E1").Value = 2.707282781600952
